# Updated symbol list on Fri Dec 30 04:56:12 UTC 2022 with GitHub Actions
#
# Applies the latest coinranking.com price/label refresh to Sheet1.
# Column D ("Price") holds numeric-looking text (stored as Text, not Number,
# in the source data) so those writes are done via a leading apostrophe
# (forces Excel to keep the cell as text) followed by resetting the cell
# style back to "Normal" so no stray NumberFormat/quote-prefix style sticks
# around on the cell. Column B/C/E values are plain text and need no
# special handling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($address, $value) {
    $cell = $ws.Range($address)
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

function Set-PlainValue($address, $value) {
    $ws.Range($address).Value = $value
}

# --- Column D "Price" refresh (numeric-looking text) -------------------
Set-TextValue "D2"  "245.58"
Set-TextValue "D3"  "23.97"
Set-TextValue "D4"  "5.248"
Set-TextValue "D5"  "0.05793"
Set-TextValue "D6"  "6.503"
Set-TextValue "D8"  "0.8172"
Set-TextValue "D9"  "0.8503"

# --- Rows 10-18: coin list shifted down by one slot ---------------------
Set-PlainValue "B10" "WazirX"
Set-PlainValue "C10" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue  "D10" "0.1365"
Set-PlainValue "E10" "9WazirXWRX"

Set-PlainValue "B11" "MandalaExchangeToken"
Set-PlainValue "C11" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue  "D11" "0.06954"
Set-PlainValue "E11" "10MandalaExchangeTokenMDX"

Set-PlainValue "B12" "LiechtensteinCryptoassetsExchange"
Set-PlainValue "C12" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue  "D12" "0.03194"
Set-PlainValue "E12" "11LiechtensteinCryptoassetsExchangeLCX"

Set-PlainValue "B13" "BitrueCoin"
Set-PlainValue "C13" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue  "D13" "0.02870"
Set-PlainValue "E13" "12BitrueCoinBTR"

Set-PlainValue "B14" "BitMartToken"
Set-PlainValue "C14" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue  "D14" "0.09377"
Set-PlainValue "E14" "13BitMartTokenBMX"

Set-PlainValue "B15" "MCDex"
Set-PlainValue "C15" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue  "D15" "3.749"
Set-PlainValue "E15" "14MCDexMCB"

Set-PlainValue "B16" "BitForexToken"
Set-PlainValue "C16" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue  "D16" "0.001524"
Set-PlainValue "E16" "15BitForexTokenBF"

Set-PlainValue "B17" "CoinExToken"
Set-PlainValue "C17" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue  "D17" "0.04722"
Set-PlainValue "E17" "16CoinExTokenCET"

Set-PlainValue "B18" "One"
Set-PlainValue "C18" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue  "D18" "0.0005958"
Set-PlainValue "E18" "17OneONE"

# --- Remaining scattered Price (and a couple of label) tweaks -----------
Set-TextValue "D19" "0.006278"
Set-TextValue "D21" "0.004546"
Set-TextValue "D22" "0.00006898"
Set-TextValue "D23" "3.501"
Set-TextValue "D24" "2.087"
Set-TextValue "D27" "0.1326"
Set-TextValue "D28" "0.0002328"
Set-TextValue "D40" "0.03654"

Set-TextValue  "D41" "0.003020"
Set-PlainValue "E41" "40KickTokenKICKWorstin24h"

Set-TextValue "D42" "0.1056"

Set-TextValue  "D43" "0.002749"
Set-PlainValue "E43" "42CEJICEJIBestin24h"

Set-TextValue "D44" "0.008056"
Set-TextValue "D45" "0.00005279"
Set-TextValue "D47" "0.3299"
Set-PlainValue "E47" "46CoinbaseStockTokenCOIN"

Set-TextValue "D48" "0.002343"
Set-TextValue "D49" "0.00002099"
Set-TextValue "D50" "0.0001999"
